# Add a "gender" column (F/M) right after the existing "Nom" (surname)
# column, and a new "SKF" value on the second row next to the room number.
#
# Before: ... G=Prenom H=Nom I=(room nr, text)              J=B203(row1)/date(row2, as L)  K=GIP2(row1)
# After:  ... G=Prenom H=Nom I=F/M (new)  J=(room nr, text)  K=B203(row1)/SKF(row2, new)    L=GIP2(row1)   M=date(row2)
#
# i.e. a brand-new, blank column is inserted at column I, shifting every
# existing column from I onward one place to the right; then the two new
# gender cells and the new "SKF" cell are populated in their new spots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank column before column I (9) - this shifts the old
# I/J/K/L columns (and their formatting / widths) one position to the right,
# turning old I -> J, old J -> K, old K -> L.
$ws.Columns.Item(9).Insert()

# Populate the newly inserted column I with the gender values.
$ws.Range("I1").Value = "F"
$ws.Range("I2").Value = "M"

# Populate the new "SKF" value that now sits in column K on row 2.
$ws.Range("K2").Value = "SKF"

# Match the author's final selection (cell K2).
[void]$ws.Range("K2").Select()
